# Auto-applied scheduled market-data refresh for the Leve profit sheets.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) per sheet
# with freshly fetched marketboard data. Values only; no formulas, formatting,
# or structural changes beyond the two cells noted below.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1612.4762
$ws.Range("I33").Value = 2021.2142
$ws.Range("J33").Value = 795
$ws.Range("K33").Value = 2021.2142
$ws.Range("L33").Value = 795
$ws.Range("M33").Value = -1792.2142
$ws.Range("N33").Value = -1253
$ws.Range("H38").Value = 268.7
$ws.Range("I38").Value = 243
$ws.Range("J38").Value = 500
$ws.Range("K38").Value = 729
$ws.Range("L38").Value = 1500
$ws.Range("M38").Value = -357
$ws.Range("N38").Value = -2244
$ws.Range("H98").Value = 3404.48
$ws.Range("I98").Value = 3657.087
$ws.Range("K98").Value = 3657.087
$ws.Range("M98").Value = -2159.087
$ws.Range("H107").Value = 17857624
$ws.Range("I107").Value = 6944811
$ws.Range("J107").Value = 83334500
$ws.Range("K107").Value = 6944811
$ws.Range("L107").Value = 83334500
$ws.Range("M107").Value = -6942891
$ws.Range("N107").Value = -83338340
$ws.Range("H111").Value = 17859728
$ws.Range("I111").Value = 25001024
$ws.Range("K111").Value = 75003072
$ws.Range("M111").Value = -75000005
$ws.Range("H122").Value = 3404.48
$ws.Range("I122").Value = 3657.087
$ws.Range("K122").Value = 10971.261
$ws.Range("M122").Value = -8521.261
$ws.Range("H135").Value = 625742.7
$ws.Range("I135").Value = 667418.9
$ws.Range("K135").Value = 6006770.100000001
$ws.Range("M135").Value = -6004235.100000001
$ws.Range("H138").Value = 4509.396
$ws.Range("I138").Value = 1040.25
$ws.Range("J138").Value = 9366.200000000001
$ws.Range("K138").Value = 3120.75
$ws.Range("L138").Value = 28098.6
$ws.Range("M138").Value = 2019.25
$ws.Range("N138").Value = -38378.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2695435.8
$ws.Range("I32").Value = 2947321.2
$ws.Range("K32").Value = 2947321.2
$ws.Range("M32").Value = -2947034.2
$ws.Range("H45").Value = 3570.158
$ws.Range("I45").Value = 1900.4166
$ws.Range("K45").Value = 1900.4166
$ws.Range("M45").Value = -1523.4166
$ws.Range("H61").Value = 4755.183
$ws.Range("I61").Value = 2519.848
$ws.Range("J61").Value = 12099.857
$ws.Range("K61").Value = 2519.848
$ws.Range("L61").Value = 12099.857
$ws.Range("M61").Value = -2307.848
$ws.Range("N61").Value = -12523.857
$ws.Range("H102").Value = 1115.8334
$ws.Range("I102").Value = 1115.8334
$ws.Range("K102").Value = 1115.8334
$ws.Range("M102").Value = 506.1666
$ws.Range("H132").Value = 3153.8315
$ws.Range("I132").Value = 1595.4058
$ws.Range("J132").Value = 8530.4
$ws.Range("K132").Value = 4786.2174
$ws.Range("L132").Value = 25591.2
$ws.Range("M132").Value = -2256.2174
$ws.Range("N132").Value = -30651.2
$ws.Range("H136").Value = 4755.183
$ws.Range("I136").Value = 2519.848
$ws.Range("J136").Value = 12099.857
$ws.Range("K136").Value = 7559.544
$ws.Range("L136").Value = 36299.571
$ws.Range("M136").Value = -5009.544
$ws.Range("N136").Value = -41399.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4329.4
$ws.Range("I105").Value = 3004.8333
$ws.Range("K105").Value = 3004.8333
$ws.Range("M105").Value = -1257.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 4080
$ws.Range("I6").Value = 4080
$ws.Range("K6").Value = 4080
$ws.Range("M6").Value = -3967
$ws.Range("H16").Value = 8390.214
$ws.Range("I16").Value = 8066.143
$ws.Range("J16").Value = 8714.286
$ws.Range("K16").Value = 8066.143
$ws.Range("L16").Value = 8714.286
$ws.Range("M16").Value = -7779.143
$ws.Range("N16").Value = -9288.286
$ws.Range("H31").Value = 15218
$ws.Range("I31").Value = 10747
$ws.Range("J31").Value = 15856.714
$ws.Range("K31").Value = 10747
$ws.Range("L31").Value = 15856.714
$ws.Range("M31").Value = -10452
$ws.Range("N31").Value = -16446.714
$ws.Range("H34").Value = 15218
$ws.Range("I34").Value = 10747
$ws.Range("J34").Value = 15856.714
$ws.Range("K34").Value = 10747
$ws.Range("L34").Value = 15856.714
$ws.Range("M34").Value = -10545
$ws.Range("N34").Value = -16260.714
$ws.Range("H58").Value = 3719.3418
$ws.Range("I58").Value = 992.0784
$ws.Range("J58").Value = 8686.857
$ws.Range("K58").Value = 992.0784
$ws.Range("L58").Value = 8686.857
$ws.Range("M58").Value = -789.0784
$ws.Range("N58").Value = -9092.857
$ws.Range("H86").Value = 5222118
$ws.Range("I86").Value = 6261841.5
$ws.Range("K86").Value = 6261841.5
$ws.Range("M86").Value = -6260718.5
$ws.Range("H89").Value = 5222118
$ws.Range("I89").Value = 6261841.5
$ws.Range("K89").Value = 31309207.5
$ws.Range("M89").Value = -31303591.5
$ws.Range("H113").Value = 8390.214
$ws.Range("I113").Value = 8066.143
$ws.Range("J113").Value = 8714.286
$ws.Range("K113").Value = 8066.143
$ws.Range("L113").Value = 8714.286
$ws.Range("M113").Value = -5896.143
$ws.Range("N113").Value = -13054.286
$ws.Range("H132").Value = 3725.6333
$ws.Range("I132").Value = 1651.4762
$ws.Range("K132").Value = 4954.4286
$ws.Range("M132").Value = -2424.4286
$ws.Range("H134").Value = 3478.7065
$ws.Range("J134").Value = 6206.5127
$ws.Range("L134").Value = 18619.5381
$ws.Range("N134").Value = -23689.5381
$ws.Range("H136").Value = 3719.3418
$ws.Range("I136").Value = 992.0784
$ws.Range("J136").Value = 8686.857
$ws.Range("K136").Value = 2976.2352
$ws.Range("L136").Value = 26060.571
$ws.Range("M136").Value = -426.2352000000001
$ws.Range("N136").Value = -31160.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2982.2104
$ws.Range("I5").Value = 2584.9524
$ws.Range("J5").Value = 3472.9412
$ws.Range("K5").Value = 7754.8572
$ws.Range("L5").Value = 10418.8236
$ws.Range("M5").Value = -7642.8572
$ws.Range("N5").Value = -10642.8236
$ws.Range("H92").Value = 5918811.5
$ws.Range("J92").Value = 5918811.5
$ws.Range("L92").Value = 17756434.5
$ws.Range("N92").Value = -17758930.5
$ws.Range("H97").Value = 312.6
$ws.Range("I97").Value = 340.75
$ws.Range("K97").Value = 1022.25
$ws.Range("M97").Value = -526.25
$ws.Range("H113").Value = 5770.4375
$ws.Range("I113").Value = 2899.8
$ws.Range("J113").Value = 7075.273
$ws.Range("K113").Value = 8699.400000000001
$ws.Range("L113").Value = 21225.819
$ws.Range("M113").Value = -6529.400000000001
$ws.Range("N113").Value = -25565.819
$ws.Range("H122").Value = 1490170.4
$ws.Range("I122").Value = 3143936
$ws.Range("J122").Value = 1781.2
$ws.Range("K122").Value = 28295424
$ws.Range("L122").Value = 16030.8
$ws.Range("M122").Value = -28292974
$ws.Range("N122").Value = -20930.8
$ws.Range("H135").Value = 2982.2104
$ws.Range("I135").Value = 2584.9524
$ws.Range("J135").Value = 3472.9412
$ws.Range("K135").Value = 23264.5716
$ws.Range("L135").Value = 31256.4708
$ws.Range("M135").Value = -20729.5716
$ws.Range("N135").Value = -36326.4708

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H102").Value = 2754.2222
$ws.Range("I102").Value = 2475.5386
$ws.Range("K102").Value = 2475.5386
$ws.Range("M102").Value = -853.5385999999999
$ws.Range("H122").Value = 45417.32
$ws.Range("I122").Value = 54357.78
$ws.Range("J122").Value = 4291.2
$ws.Range("K122").Value = 163073.34
$ws.Range("L122").Value = 12873.6
$ws.Range("M122").Value = -160623.34
$ws.Range("N122").Value = -17773.6
$ws.Range("H126").Value = 5784.7144
$ws.Range("I126").Value = 4331.8887
$ws.Range("J126").Value = 8399.799999999999
$ws.Range("K126").Value = 12995.6661
$ws.Range("L126").Value = 25199.4
$ws.Range("M126").Value = -10525.6661
$ws.Range("N126").Value = -30139.4
$ws.Range("H134").Value = 77264.5
$ws.Range("J134").Value = 77264.5
$ws.Range("L134").Value = 231793.5
$ws.Range("N134").Value = -236863.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4762.9814
$ws.Range("I7").Value = 3673.2144
$ws.Range("K7").Value = 3673.2144
$ws.Range("M7").Value = -3561.2144
$ws.Range("H68").Value = 7013.8945
$ws.Range("I68").Value = 5825.2856
$ws.Range("J68").Value = 7707.25
$ws.Range("K68").Value = 5825.2856
$ws.Range("L68").Value = 7707.25
$ws.Range("M68").Value = -5076.2856
$ws.Range("N68").Value = -9205.25
$ws.Range("H71").Value = 7013.8945
$ws.Range("I71").Value = 5825.2856
$ws.Range("J71").Value = 7707.25
$ws.Range("K71").Value = 29126.428
$ws.Range("L71").Value = 38536.25
$ws.Range("M71").Value = -25382.428
$ws.Range("N71").Value = -46024.25
$ws.Range("H126").Value = 4762.9814
$ws.Range("I126").Value = 3673.2144
$ws.Range("K126").Value = 11019.6432
$ws.Range("M126").Value = -8549.643199999999
$ws.Range("H136").Value = 5446.6616
$ws.Range("I136").Value = 1717.6511
$ws.Range("K136").Value = 5152.9533
$ws.Range("M136").Value = -2602.9533

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4120
$ws.Range("I96").Value = 1600
$ws.Range("J96").Value = 4750
$ws.Range("K96").Value = 1600
$ws.Range("L96").Value = 4750
$ws.Range("M96").Value = -227
$ws.Range("N96").Value = -7496
$ws.Range("H107").Value = 7937235.5
$ws.Range("I107").Value = 376.89474
$ws.Range("K107").Value = 1130.68422
$ws.Range("M107").Value = 789.3157799999999
$ws.Range("H132").Value = 19611396
$ws.Range("I132").Value = 22225448
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 66676344
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -66673814
$ws.Range("N132").Value = -23057
$ws.Range("H136").Value = 15306042
$ws.Range("I136").Value = 23256704
$ws.Range("J136").Value = 441761.25
$ws.Range("K136").Value = 69770112
$ws.Range("L136").Value = 1325283.75
$ws.Range("M136").Value = -69767562
$ws.Range("N136").Value = -1330383.75

Write-Output "Sheets updated via scheduled runner."
